# --- Apply the COVID-19 country-data refresh (matches the source commit "Update countries & provincias Spain") ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last refreshed" timestamp label in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 06:22"

# 2) Replace the country table (A4:H215) with the refreshed figures, already
#    re-sorted by "Casos totales" (column B) descending -- this reproduces both
#    the updated case counts and the row reshuffle caused by the new ranking.
#    Columns: Pais, Casos totales, Nuevos casos, Casos activos, Recuperados,
#             Casos criticos, Muertes hoy, Muertes
$rows = @(
  ,@("Estados Unidos", 400540, 205, 21711, 365972, 9169, 16, 12857)
  ,@("España", 141942, 0, 43208, 84689, 7069, 0, 14045)
  ,@("Italia", 135586, 0, 24392, 94067, 3792, 0, 17127)
  ,@("Francia", 109069, 0, 19337, 79404, 7131, 0, 10328)
  ,@("Alemania", 107663, 0, 36081, 69566, 4895, 0, 2016)
  ,@("China", 81802, 62, 77279, 1190, 189, 2, 3333)
  ,@("Iran", 62589, 0, 27039, 31678, 3987, 0, 3872)
  ,@("Reino Unido", 55242, 0, 135, 48948, 1559, 0, 6159)
  ,@("Turquia", 34109, 0, 1582, 31802, 1474, 0, 725)
  ,@("Suiza", 22253, 0, 8704, 12728, 391, 0, 821)
  ,@("Belgica", 22194, 0, 4157, 16002, 1260, 0, 2035)
  ,@("Paises Bajos", 19580, 0, 250, 17229, 1424, 0, 2101)
  ,@("Canada", 17897, 0, 4028, 13488, 426, 0, 381)
  ,@("Brasil", 14049, 15, 127, 13234, 296, 2, 688)
  ,@("Austria", 12639, 0, 4046, 8350, 243, 0, 243)
  ,@("Portugal", 12442, 0, 184, 11913, 271, 0, 345)
  ,@("Corea del Sur", 10384, 53, 6776, 3408, 55, 8, 200)
  ,@("Israel", 9248, 0, 770, 8413, 149, 0, 65)
  ,@("Suecia", 7693, 0, 205, 6897, 640, 0, 591)
  ,@("Rusia", 7497, 0, 494, 6945, 8, 0, 58)
  ,@("Noruega", 6086, 0, 32, 5965, 78, 0, 89)
  ,@("Australia", 5988, 0, 2547, 3392, 93, 0, 49)
  ,@("Irlanda", 5709, 0, 25, 5474, 165, 0, 210)
  ,@("India", 5356, 5, 468, 4728, 0, 0, 160)
  ,@("Chile", 5116, 0, 898, 4175, 337, 0, 43)
  ,@("Dinamarca", 5071, 0, 1491, 3377, 127, 0, 203)
  ,@("Chequia", 5017, 0, 172, 4757, 86, 0, 88)
  ,@("Polonia", 4848, 0, 191, 4528, 50, 0, 129)
  ,@("Rumania", 4417, 0, 460, 3760, 274, 0, 197)
  ,@("Japon", 4257, 0, 622, 3542, 80, 0, 93)
  ,@("Pakistan", 4035, 0, 429, 3549, 28, 0, 57)
  ,@("Ecuador", 3995, 0, 140, 3635, 156, 0, 220)
  ,@("Malasia", 3963, 0, 1321, 2579, 92, 0, 63)
  ,@("Filipinas", 3764, 0, 84, 3503, 1, 0, 177)
  ,@("Luxemburgo", 2970, 0, 500, 2426, 35, 0, 44)
  ,@("Peru", 2954, 0, 1301, 1546, 109, 0, 107)
  ,@("Arabia Saudita", 2795, 0, 615, 2139, 41, 0, 41)
  ,@("Mexico", 2785, 346, 633, 2011, 89, 16, 141)
  ,@("Indonesia", 2738, 0, 204, 2313, 0, 0, 221)
  ,@("Serbia", 2447, 0, 118, 2268, 109, 0, 61)
  ,@("Emiratos Arabes Unidos", 2359, 0, 186, 2161, 1, 0, 12)
  ,@("Finlandia", 2308, 0, 300, 1974, 81, 0, 34)
  ,@("Tailandia", 2258, 0, 888, 1343, 61, 0, 27)
  ,@("Panama", 2249, 149, 16, 2174, 91, 4, 59)
  ,@("Catar", 2057, 0, 150, 1901, 37, 0, 6)
  ,@("Republica Dominicana", 1956, 0, 36, 1822, 147, 0, 98)
  ,@("Grecia", 1832, 0, 269, 1482, 90, 0, 81)
  ,@("Colombia", 1780, 0, 100, 1630, 76, 0, 50)
  ,@("Sudafrica", 1749, 0, 95, 1641, 7, 0, 13)
  ,@("Argentina", 1715, 0, 338, 1317, 96, 0, 60)
  ,@("Islandia", 1586, 0, 559, 1021, 11, 0, 6)
  ,@("Singapur", 1481, 0, 377, 1098, 29, 0, 6)
  ,@("Argelia", 1468, 0, 113, 1162, 46, 0, 193)
  ,@("Ucrania", 1462, 0, 28, 1389, 16, 0, 45)
  ,@("Egipto", 1450, 0, 276, 1080, 0, 0, 94)
  ,@("Croacia", 1282, 0, 167, 1097, 35, 0, 18)
  ,@("Nueva Zelanda", 1210, 50, 282, 927, 4, 0, 1)
  ,@("Marruecos", 1184, 0, 93, 1001, 1, 0, 90)
  ,@("Estonia", 1149, 0, 69, 1059, 12, 0, 21)
  ,@("Irak", 1122, 0, 373, 684, 0, 0, 65)
  ,@("Eslovenia", 1059, 0, 102, 921, 30, 0, 36)
  ,@("Moldavia", 1056, 0, 40, 994, 80, 0, 22)
  ,@("Hong Kong", 936, 0, 236, 696, 12, 0, 4)
  ,@("Hungria", 895, 78, 94, 754, 17, 0, 47)
  ,@("Lituania", 880, 0, 8, 857, 11, 0, 15)
  ,@("Bielorrusia", 861, 0, 54, 794, 31, 0, 13)
  ,@("Armenia", 853, 0, 87, 758, 30, 0, 8)
  ,@("Barein", 811, 0, 458, 348, 4, 0, 5)
  ,@("Bosnia y Herzegovina", 764, 0, 68, 663, 4, 0, 33)
  ,@("Kuwait", 743, 0, 105, 637, 23, 0, 1)
  ,@("Azerbaiyan", 717, 0, 44, 665, 23, 0, 8)
  ,@("Crucero", 712, 0, 619, 82, 10, 0, 11)
  ,@("Kazajistan", 697, 0, 51, 640, 16, 0, 6)
  ,@("Camerun", 685, 0, 60, 616, 0, 0, 9)
  ,@("Tunez", 623, 0, 25, 575, 39, 0, 23)
  ,@("Republica de Macedonia", 599, 0, 30, 543, 15, 0, 26)
  ,@("Eslovaquia", 581, 0, 13, 566, 3, 0, 2)
  ,@("Bulgaria", 577, 0, 42, 512, 21, 0, 23)
  ,@("Letonia", 548, 0, 16, 530, 5, 0, 2)
  ,@("Libano", 548, 0, 62, 467, 29, 0, 19)
  ,@("Principado de Andorra", 545, 0, 39, 484, 17, 0, 22)
  ,@("Uzbekistan", 520, 0, 30, 488, 8, 0, 2)
  ,@("Republica de Chipre", 494, 0, 47, 438, 13, 0, 9)
  ,@("Costa Rica", 483, 0, 24, 457, 14, 0, 2)
  ,@("Uruguay", 424, 0, 150, 267, 14, 0, 7)
  ,@("Afganistan", 423, 0, 18, 391, 0, 0, 14)
  ,@("Cuba", 396, 0, 27, 358, 15, 0, 11)
  ,@("Burkina Faso", 384, 0, 127, 238, 0, 0, 19)
  ,@("Albania", 383, 0, 131, 230, 7, 0, 22)
  ,@("Taiwan", 376, 0, 61, 310, 0, 0, 5)
  ,@("Oman", 371, 0, 67, 302, 3, 0, 2)
  ,@("Reunion", 358, 0, 40, 318, 4, 0, 0)
  ,@("Jordania", 353, 0, 138, 209, 5, 0, 6)
  ,@("Costa de Marfil", 349, 0, 41, 305, 0, 0, 3)
  ,@("Honduras", 312, 7, 6, 284, 10, 0, 22)
  ,@("Malta", 293, 0, 5, 288, 4, 0, 0)
  ,@("Ghana", 287, 0, 31, 251, 2, 0, 5)
  ,@("San Marino", 279, 0, 40, 205, 14, 0, 34)
  ,@("Niger", 278, 0, 26, 241, 0, 0, 11)
  ,@("Mauricio", 268, 0, 8, 253, 3, 0, 7)
  ,@("Estado de Palestina", 261, 0, 42, 218, 0, 0, 1)
  ,@("Nigeria", 254, 0, 44, 204, 2, 0, 6)
  ,@("Vietnam", 251, 0, 122, 129, 8, 0, 0)
  ,@("Montenegro", 241, 0, 4, 235, 7, 0, 2)
  ,@("Senegal", 237, 0, 105, 130, 1, 0, 2)
  ,@("Kirguistan", 228, 0, 33, 191, 5, 0, 4)
  ,@("Bolivia", 210, 16, 2, 193, 3, 1, 15)
  ,@("Georgia", 196, 0, 46, 147, 6, 0, 3)
  ,@("Sri Lanka", 185, 0, 42, 137, 5, 0, 6)
  ,@("Islas Feroe", 184, 0, 120, 64, 2, 0, 0)
  ,@("Consejo Danes para los Refugiados", 180, 0, 9, 153, 0, 0, 18)
  ,@("Kenia", 172, 0, 7, 159, 2, 0, 6)
  ,@("Mayotte", 171, 0, 22, 147, 3, 0, 2)
  ,@("Venezuela", 166, 0, 65, 94, 6, 0, 7)
  ,@("Banglades", 164, 0, 33, 114, 1, 0, 17)
  ,@("Martinica", 152, 0, 50, 98, 20, 0, 4)
  ,@("Isla de Man", 150, 0, 73, 76, 6, 0, 1)
  ,@("Guinea", 144, 0, 5, 139, 0, 0, 0)
  ,@("Guadalupe", 139, 0, 31, 101, 14, 0, 7)
  ,@("Brunei", 135, 0, 85, 49, 3, 0, 1)
  ,@("Paraguay", 119, 4, 15, 99, 1, 0, 5)
  ,@("Camboya", 117, 2, 63, 54, 1, 0, 0)
  ,@("Gibraltar", 113, 0, 60, 53, 0, 0, 0)
  ,@("Trinidad yTobago", 107, 0, 1, 98, 0, 0, 8)
  ,@("Ruanda", 105, 0, 7, 98, 0, 0, 0)
  ,@("Republica de Yibuti", 90, 0, 9, 81, 0, 0, 0)
  ,@("Madagascar", 88, 0, 7, 81, 6, 0, 0)
  ,@("Guatemala", 80, 3, 17, 60, 3, 0, 3)
  ,@("Monaco", 79, 0, 4, 74, 4, 0, 1)
  ,@("El Salvador", 78, 0, 5, 69, 4, 0, 4)
  ,@("Liechtenstein", 78, 0, 55, 22, 0, 0, 1)
  ,@("Guayana Francesa", 77, 0, 34, 43, 1, 0, 0)
  ,@("Aruba", 74, 0, 14, 60, 0, 0, 0)
  ,@("Togo", 65, 0, 23, 39, 0, 0, 3)
  ,@("Barbados", 63, 0, 6, 54, 4, 0, 3)
  ,@("Jamaica", 63, 0, 9, 51, 0, 0, 3)
  ,@("Mali", 56, 0, 12, 39, 0, 0, 5)
  ,@("Uganda", 52, 0, 0, 52, 0, 0, 0)
  ,@("Etiopia", 52, 0, 4, 46, 1, 0, 2)
  ,@("Polinesia Francesa", 47, 0, 0, 47, 0, 0, 0)
  ,@("Islas Caimanes", 45, 0, 6, 38, 0, 0, 1)
  ,@("Congo", 45, 0, 2, 38, 0, 0, 5)
  ,@("Macao", 44, 0, 10, 34, 1, 0, 0)
  ,@("San Martin (Parte Holandesa)", 40, 0, 1, 33, 2, 0, 6)
  ,@("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
  ,@("Zambia", 39, 0, 7, 31, 0, 0, 1)
  ,@("Bermudas", 39, 0, 17, 20, 0, 0, 2)
  ,@("Bahamas", 36, 3, 5, 25, 1, 0, 6)
  ,@("Guinea-Bisau", 33, 0, 0, 33, 0, 0, 0)
  ,@("Guyana", 33, 0, 8, 20, 4, 0, 5)
  ,@("Guam", 32, 0, 0, 31, 0, 0, 1)
  ,@("San Martin (Parte Francesa)", 32, 0, 7, 23, 6, 0, 2)
  ,@("Eritrea", 31, 0, 0, 31, 0, 0, 0)
  ,@("Gabon", 30, 0, 1, 28, 0, 0, 1)
  ,@("Benin", 26, 0, 5, 20, 0, 0, 1)
  ,@("Haiti", 25, 0, 0, 24, 0, 0, 1)
  ,@("Tanzania", 24, 0, 5, 18, 0, 0, 1)
  ,@("Birmania", 22, 0, 0, 21, 0, 0, 1)
  ,@("Libia", 20, 0, 1, 18, 0, 0, 1)
  ,@("Antigua y Barbuda", 19, 0, 0, 18, 1, 0, 1)
  ,@("Siria", 19, 0, 3, 14, 0, 0, 2)
  ,@("Maldivas", 19, 0, 13, 6, 0, 0, 0)
  ,@("Nueva Caledonia", 18, 0, 1, 17, 0, 0, 0)
  ,@("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
  ,@("Angola", 17, 0, 2, 13, 0, 0, 2)
  ,@("Namibia", 16, 0, 3, 13, 0, 0, 0)
  ,@("Guinea Ecuatorial", 16, 0, 3, 13, 0, 0, 0)
  ,@("Fiyi", 15, 0, 0, 15, 0, 0, 0)
  ,@("Dominica", 15, 0, 1, 14, 0, 0, 0)
  ,@("Mongolia", 15, 0, 4, 11, 0, 0, 0)
  ,@("Laos", 14, 0, 0, 14, 0, 0, 0)
  ,@("Santa Lucia", 14, 0, 1, 13, 0, 0, 0)
  ,@("Sudan", 14, 0, 2, 10, 0, 0, 2)
  ,@("Liberia", 14, 0, 3, 8, 0, 0, 3)
  ,@("Curazao", 13, 0, 7, 5, 0, 0, 1)
  ,@("Granada", 12, 0, 0, 12, 2, 0, 0)
  ,@("San Cristobal y Nieves", 11, 0, 0, 11, 0, 0, 0)
  ,@("Seychelles", 11, 0, 0, 11, 0, 0, 0)
  ,@("Zimbabue", 11, 0, 0, 9, 0, 0, 2)
  ,@("Groenlandia", 11, 0, 10, 1, 0, 0, 0)
  ,@("Mozambique", 10, 0, 1, 9, 0, 0, 0)
  ,@("Republica del Chad", 10, 0, 2, 8, 0, 0, 0)
  ,@("Surinam", 10, 0, 3, 6, 0, 0, 1)
  ,@("Suazilandia", 10, 0, 4, 6, 0, 0, 0)
  ,@("Nepal", 9, 0, 1, 8, 0, 0, 0)
  ,@("Montserrat", 9, 0, 0, 7, 0, 0, 2)
  ,@("Republica de Africa Central", 8, 0, 0, 8, 0, 0, 0)
  ,@("Islas Turcas y Caicos", 8, 0, 0, 7, 0, 0, 1)
  ,@("San Vicente y las Granadinas", 8, 0, 1, 7, 0, 0, 0)
  ,@("Malaui", 8, 0, 0, 7, 1, 0, 1)
  ,@("Somalia", 8, 0, 1, 7, 0, 0, 0)
  ,@("Santa Sede", 7, 0, 0, 7, 0, 0, 0)
  ,@("Belice", 7, 0, 0, 6, 1, 0, 1)
  ,@("Cabo Verde", 7, 0, 1, 5, 0, 0, 1)
  ,@("Sierra Leona", 6, 0, 0, 6, 0, 0, 0)
  ,@("Nicaragua", 6, 0, 0, 5, 0, 0, 1)
  ,@("Botsuana", 6, 0, 0, 5, 0, 0, 1)
  ,@("San Bartolome", 6, 0, 1, 5, 0, 0, 0)
  ,@("Mauritania", 6, 0, 2, 3, 0, 0, 1)
  ,@("Butan", 5, 0, 2, 3, 0, 0, 0)
  ,@("Sahara Occidental", 4, 0, 0, 4, 0, 0, 0)
  ,@("Santo Tome y Principe", 4, 0, 0, 4, 0, 0, 0)
  ,@("Gambia", 4, 0, 2, 1, 0, 0, 1)
  ,@("Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0)
  ,@("Anguila", 3, 0, 0, 3, 0, 0, 0)
  ,@("Burundi", 3, 0, 0, 3, 0, 0, 0)
  ,@("Sudan del Sur", 2, 0, 0, 2, 0, 0, 0)
  ,@("Bonaire, San Eustaquio y Saba", 2, 0, 0, 2, 0, 0, 0)
  ,@("Islas Malvinas", 2, 0, 0, 2, 0, 0, 0)
  ,@("Papua Nueva Guinea", 2, 0, 0, 2, 0, 0, 0)
  ,@("Timor Oriental", 1, 0, 0, 1, 0, 0, 0)
  ,@("San Pedro y Miquelon", 1, 0, 0, 1, 0, 0, 0)
)

$numRows = $rows.Count
$numCols = 8
$data = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
  for ($j = 0; $j -lt $numCols; $j++) {
    $data[$i, $j] = $rows[$i][$j]
  }
}

$ws.Range("A4:H215").Value = $data

